$d = $word.ActiveDocument

$replacements = @(
    @("48÷4=", "68÷2="),
    @("41÷4=", "32÷7="),
    @("94÷7=", "32÷2="),
    @("54÷5=", "97÷6="),
    @("42÷4=", "95÷4="),
    @("64÷4=", "68÷4="),
    @("45÷9=", "63÷5="),
    @("77÷7=", "25÷3="),
    @("69÷7=", "95÷8="),
    @("79÷7=", "62÷5="),
    @("97÷7=", "84÷5="),
    @("72÷9=", "98÷2="),
    @("26÷2=", "18÷3="),
    @("94÷8=", "84÷8="),
    @("96÷8=", "88÷9="),
    @("87÷4=", "60÷9="),
    @("57÷9=", "76÷9="),
    @("22÷9=", "74÷6="),
    @("49÷6=", "85÷3="),
    @("45÷7=", "87÷2="),
    @("30÷7=", "11÷9="),
    @("81÷9=", "65÷9="),
    @("40÷2=", "42÷3="),
    @("19÷5=", "69÷5="),
    @("28÷7=", "75÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
